$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows 42-82 of sequential-extraction data (Fe and Mn parameters)
$ws.Range("A42").Value = 1
$ws.Range("B42").Value = "30-35"
$ws.Range("C42").Value = 1.5083
$ws.Range("D42").Value = "H2O"
$ws.Range("E42").Value = "Fe"
$ws.Range("F42").Value = 0.59

$ws.Range("A43").Value = 2
$ws.Range("B43").Value = "30-35"
$ws.Range("C43").Value = 1.3788
$ws.Range("D43").Value = "H2O"
$ws.Range("E43").Value = "Fe"
$ws.Range("F43").Value = 0.55

$ws.Range("A44").Value = 3
$ws.Range("B44").Value = "30-35"
$ws.Range("C44").Value = 1.4689
$ws.Range("D44").Value = "H2O"
$ws.Range("E44").Value = "Fe"
$ws.Range("F44").Value = 0.55

$ws.Range("A45").Value = 4
$ws.Range("B45").Value = "30-35"
$ws.Range("C45").Value = 1.548
$ws.Range("D45").Value = "H2O"
$ws.Range("E45").Value = "Fe"
$ws.Range("F45").Value = 0.51

$ws.Range("A46").Value = "blank"
$ws.Range("D46").Value = "H2O"
$ws.Range("E46").Value = "Fe"
$ws.Range("F46").Value = 0

$ws.Range("A47").Value = "blank"
$ws.Range("D47").Value = "H2O"
$ws.Range("E47").Value = "Fe"
$ws.Range("F47").Value = 0

$ws.Range("A48").Value = 1
$ws.Range("B48").Value = "30-35"
$ws.Range("C48").Value = 1.5083
$ws.Range("D48").Value = "Bipy"
$ws.Range("E48").Value = "Fe"
$ws.Range("F48").Value = 7.29

$ws.Range("A49").Value = 2
$ws.Range("B49").Value = "30-35"
$ws.Range("C49").Value = 1.3788
$ws.Range("D49").Value = "Bipy"
$ws.Range("E49").Value = "Fe"
$ws.Range("F49").Value = 6.15

$ws.Range("A50").Value = 3
$ws.Range("B50").Value = "30-35"
$ws.Range("C50").Value = 1.4689
$ws.Range("D50").Value = "Bipy"
$ws.Range("E50").Value = "Fe"
$ws.Range("F50").Value = 6.95

$ws.Range("A51").Value = 4
$ws.Range("B51").Value = "30-35"
$ws.Range("C51").Value = 1.548
$ws.Range("D51").Value = "Bipy"
$ws.Range("E51").Value = "Fe"
$ws.Range("F51").Value = 6.19

$ws.Range("A52").Value = "blank"
$ws.Range("D52").Value = "Bipy"
$ws.Range("E52").Value = "Fe"
$ws.Range("F52").Value = 0

$ws.Range("A53").Value = "blank"
$ws.Range("D53").Value = "Bipy"
$ws.Range("E53").Value = "Fe"
$ws.Range("F53").Value = 0

$ws.Range("A54").Value = 1
$ws.Range("B54").Value = "30-35"
$ws.Range("C54").Value = 1.5083
$ws.Range("D54").Value = "BD"
$ws.Range("E54").Value = "Fe"
$ws.Range("F54").Value = 27.8

$ws.Range("A55").Value = 2
$ws.Range("B55").Value = "30-35"
$ws.Range("C55").Value = 1.3788
$ws.Range("D55").Value = "BD"
$ws.Range("E55").Value = "Fe"
$ws.Range("F55").Value = 23.9

$ws.Range("A56").Value = 3
$ws.Range("B56").Value = "30-35"
$ws.Range("C56").Value = 1.4689
$ws.Range("D56").Value = "BD"
$ws.Range("E56").Value = "Fe"
$ws.Range("F56").Value = 34.2

$ws.Range("A57").Value = 4
$ws.Range("B57").Value = "30-35"
$ws.Range("C57").Value = 1.548
$ws.Range("D57").Value = "BD"
$ws.Range("E57").Value = "Fe"
$ws.Range("F57").Value = 18.9

$ws.Range("A58").Value = "blank"
$ws.Range("D58").Value = "BD"
$ws.Range("E58").Value = "Fe"
$ws.Range("F58").Value = 0

$ws.Range("A59").Value = "blank"
$ws.Range("D59").Value = "BD"
$ws.Range("E59").Value = "Fe"
$ws.Range("F59").Value = 0

$ws.Range("A60").Value = 1
$ws.Range("B60").Value = "30-35"
$ws.Range("C60").Value = 1.5083
$ws.Range("D60").Value = "H2O"
$ws.Range("E60").Value = "Mn"
$ws.Range("F60").Value = 0

$ws.Range("A61").Value = 2
$ws.Range("B61").Value = "30-35"
$ws.Range("C61").Value = 1.3788
$ws.Range("D61").Value = "H2O"
$ws.Range("E61").Value = "Mn"
$ws.Range("F61").Value = 0

$ws.Range("A62").Value = 3
$ws.Range("B62").Value = "30-35"
$ws.Range("C62").Value = 1.4689
$ws.Range("D62").Value = "H2O"
$ws.Range("E62").Value = "Mn"
$ws.Range("F62").Value = 0

$ws.Range("A63").Value = 4
$ws.Range("B63").Value = "30-35"
$ws.Range("C63").Value = 1.548
$ws.Range("D63").Value = "H2O"
$ws.Range("E63").Value = "Mn"
$ws.Range("F63").Value = 0

$ws.Range("A64").Value = "blank"
$ws.Range("D64").Value = "H2O"
$ws.Range("E64").Value = "Mn"
$ws.Range("F64").Value = 0

$ws.Range("A65").Value = "blank"
$ws.Range("D65").Value = "H2O"
$ws.Range("E65").Value = "Mn"
$ws.Range("F65").Value = 0

$ws.Range("A66").Value = 1
$ws.Range("B66").Value = "30-35"
$ws.Range("C66").Value = 1.5083
$ws.Range("D66").Value = "Bipy"
$ws.Range("E66").Value = "Mn"
$ws.Range("F66").Value = 0.19

$ws.Range("A67").Value = 2
$ws.Range("B67").Value = "30-35"
$ws.Range("C67").Value = 1.3788
$ws.Range("D67").Value = "Bipy"
$ws.Range("E67").Value = "Mn"
$ws.Range("F67").Value = 0.2

$ws.Range("A68").Value = 3
$ws.Range("B68").Value = "30-35"
$ws.Range("C68").Value = 1.4689
$ws.Range("D68").Value = "Bipy"
$ws.Range("E68").Value = "Mn"
$ws.Range("F68").Value = 0.2

$ws.Range("A69").Value = 4
$ws.Range("B69").Value = "30-35"
$ws.Range("C69").Value = 1.548
$ws.Range("D69").Value = "Bipy"
$ws.Range("E69").Value = "Mn"
$ws.Range("F69").Value = 0.18

$ws.Range("A70").Value = "blank"
$ws.Range("D70").Value = "Bipy"
$ws.Range("E70").Value = "Mn"
$ws.Range("F70").Value = 0

$ws.Range("A71").Value = "blank"
$ws.Range("D71").Value = "Bipy"
$ws.Range("E71").Value = "Mn"
$ws.Range("F71").Value = 0

$ws.Range("A72").Value = 1
$ws.Range("B72").Value = "30-35"
$ws.Range("C72").Value = 1.5083
$ws.Range("D72").Value = "BD"
$ws.Range("E72").Value = "Mn"
$ws.Range("F72").Value = 3.81

$ws.Range("A73").Value = 2
$ws.Range("B73").Value = "30-35"
$ws.Range("C73").Value = 1.3788
$ws.Range("D73").Value = "BD"
$ws.Range("E73").Value = "Mn"
$ws.Range("F73").Value = 3.41

$ws.Range("A74").Value = 3
$ws.Range("B74").Value = "30-35"
$ws.Range("C74").Value = 1.4689
$ws.Range("D74").Value = "BD"
$ws.Range("E74").Value = "Mn"
$ws.Range("F74").Value = 4.54

$ws.Range("A75").Value = 4
$ws.Range("B75").Value = "30-35"
$ws.Range("C75").Value = 1.548
$ws.Range("D75").Value = "BD"
$ws.Range("E75").Value = "Mn"
$ws.Range("F75").Value = 3.01

$ws.Range("A76").Value = "blank"
$ws.Range("D76").Value = "BD"
$ws.Range("E76").Value = "Mn"
$ws.Range("F76").Value = 0

$ws.Range("A77").Value = "blank"
$ws.Range("D77").Value = "BD"
$ws.Range("E77").Value = "Mn"
$ws.Range("F77").Value = 0

$ws.Range("A78").Value = 1
$ws.Range("B78").Value = "30-35"
$ws.Range("C78").Value = 1.5083
$ws.Range("D78").Value = "NaOH"
$ws.Range("F78").Value = 0.95

$ws.Range("A79").Value = 2
$ws.Range("B79").Value = "30-35"
$ws.Range("C79").Value = 1.3788
$ws.Range("D79").Value = "NaOH"
$ws.Range("F79").Value = 0.94

$ws.Range("A80").Value = 5
$ws.Range("B80").Value = "30-35"
$ws.Range("C80").Value = 3.0169
$ws.Range("D80").Value = "NaOH"
$ws.Range("F80").Value = 0.75

$ws.Range("A81").Value = "blank"
$ws.Range("D81").Value = "NaOH"
$ws.Range("F81").Value = 0

$ws.Range("A82").Value = "blank"
$ws.Range("D82").Value = "NaOH"
$ws.Range("F82").Value = 0

# Autofit column A to match bestFit width
$ws.Columns.Item(1).AutoFit() | Out-Null

# Update selection to F75 (also clears the old topLeftCell scroll position)
$ws.Range("F75").Select()